$d = $word.ActiveDocument

# The document contains a run of three consecutive empty paragraphs
# (justified "both", 26pt/sz=52) right after the "X(6) = 0.2, X(7) = 0.2,
# X(8) = 0.2" paragraph and before the final (big, red, bold) heading
# paragraph. The edit removes two of those three empty paragraphs,
# leaving only a single blank paragraph in their place.

# Find the paragraph that contains the "X(8) = 0.2" text to anchor the
# location robustly instead of relying on fixed indices.
$target = $null
For ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*X(8) = 0.2*") {
        $target = $i
    }
}

# The two empty paragraphs immediately following the anchor paragraph
# (index $target + 2 and $target + 3, i.e. the second and third blank
# paragraphs of the run of three) are deleted, keeping the first blank
# paragraph (index $target + 1) intact. Delete from the highest index
# down so earlier indices remain valid.
$d.Paragraphs.Item($target + 3).Range.Delete()
$d.Paragraphs.Item($target + 2).Range.Delete()
